$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A9 already has style s=3 (applyNumberFormat=1, numFmtId=0) pre-edit
$ws.Range("A9").Font.Underline = 2   # xlUnderlineStyleSingle
Write-Output "done"
